# Add a new "is_dominant" row to the Scoring file schema's "Columns" sheet,
# right after the "is_interaction" row (i.e. as new row 12), pushing every
# row below it down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at position 12 - Excel copies the formatting of the
# row above (row 11, "is_interaction") onto the new blank row, and shifts
# every row from 12 downward (and the mergeCells below) down by one.
$ws.Rows("12:12").Insert()

# Populate the three schema columns for the new "is_dominant" flag.
$ws.Range("A12").Value = "is_dominant"
$ws.Range("B12").Value = "FLAG: Dominant Inheritance Model"
$ws.Range("C12").Value = "This is a TRUE/FALSE variable that flags whether the weight should be added to the PGS sum if there is at least 1 copy of the effect allele (e.g. it is a dominant allele)."

# The field-description cell uses a slightly different dark-grey Arial run
# (as if pasted in from elsewhere), matching the other description cells'
# size/typeface but with its own font color.
$ws.Range("C12").Font.Name = "Arial"
$ws.Range("C12").Font.Size = 12
$ws.Range("C12").Font.Color = 1907741

# Leave the selection on the newly added description cell.
$ws.Range("C12").Select() | Out-Null
